# Update month_data/August.xlsx:
#  - P1 header text changes from "06" to "06.7" (kept as text, not auto-converted to a number)
#  - New column Q is added:
#      Q1 header = "18" (kept as text, same header style as the rest of row 1)
#      Q2:Q22    = attendance flags (numbers)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- P1: "06" -> "06.7" -------------------------------------------------
# A plain assignment of a numeric-looking string ("06.7") would be auto-converted
# to the number 6.7 by Excel. Route the text through a scratch cell using a
# leading apostrophe (forces text), copy/paste the *value only* into P1, then
# restore P1's original header formatting (border/bold/center) by pasting the
# format from the neighboring header cell O1.
$ws.Range("BZ1").Value = "'06.7"
$ws.Range("BZ1").Copy()
$ws.Range("P1").PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = 0
$ws.Range("O1").Copy()
$ws.Range("P1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Q1: new header "18" ------------------------------------------------
$ws.Range("BZ1").Value = "'18"
$ws.Range("BZ1").Copy()
$ws.Range("Q1").PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = 0
$ws.Range("O1").Copy()
$ws.Range("Q1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Clean up the scratch cell so it doesn't affect the sheet's used range.
$ws.Range("BZ1").Clear()

# --- Q2:Q22 data for the new day ---------------------------------------
$values = @(1,1,1,1,1,0,1,1,1,1,1,1,1,1,1,1,1,1,1,1,0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 17).Value = $values[$i]
}
